$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new working-hours entry (2014-04-16, 14:15-15:30) is being inserted
# right before the blank spacer / summary block, pushing the spacer and
# the three summary rows (sum [min], sum [h], sum [working weeks]) down
# by one row.
$ws.Rows.Item(130).Insert()

$ws.Range("A130").Value = 2014
$ws.Range("B130").Value = 4
$ws.Range("C130").Value = 16
$ws.Range("D130").Value = 0.59375
$ws.Range("E130").Value = 0.64583333333333337
$ws.Range("F130").Formula = "=(E130-D130)*24*60"
$ws.Range("G130").Formula = "=F130/60"

# The grand-total SUM, which used to stop at the (now shifted) blank
# spacer row 130, needs to keep stopping at the blank spacer row, which
# is now row 131.
$ws.Range("F132").Formula = "=SUM(F2:F131)"

$ws.Range("F130").Select()

$wb.Save()
